$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("C2").Value = 40.35273080373261
$ws.Range("E2").Value = 0.05768901818751787
$ws.Range("F2").Value = 39.7519845137213
$ws.Range("G2").Value = 39.17885991166446
$ws.Range("H2").Value = 40.33773556790973
$ws.Range("I2").Value = 0.0008155990910973818
$ws.Range("J2").Value = 0.0007160354820433255
$ws.Range("K2").Value = 0.001006965825962904
$ws.Range("L2").Value = 0.05773987243909963
$ws.Range("M2").Value = 0.05732369019961255
$ws.Range("N2").Value = 0.0581650361216453

# Update row 3 values
$ws.Range("F3").Value = 0.102597563863823
$ws.Range("G3").Value = 0.01127190191069068
$ws.Range("H3").Value = 0.2184275124632387
$ws.Range("I3").Value = 0.09367986748464503
$ws.Range("J3").Value = 0.01026609838128846
$ws.Range("K3").Value = 0.1995491589942272
$ws.Range("L3").Value = 0.1098090001372114
$ws.Range("M3").Value = 0.01210427335757139
$ws.Range("N3").Value = 0.2336030159744603

# Add new row 4
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Beta + Gamma"
$ws.Range("C4").Value = 40.35273080373261
$ws.Range("D4").Value = 0.0007432820064133916
$ws.Range("E4").Value = 0.05768901818751787
$ws.Range("F4").Value = 39.85458207758511
$ws.Range("G4").Value = 39.19013181357516
$ws.Range("H4").Value = 40.55616308037295
$ws.Range("I4").Value = 0.09449546657574241
$ws.Range("J4").Value = 0.01098213386333178
$ws.Range("K4").Value = 0.2005561248201901
$ws.Range("L4").Value = 0.167548872576311
$ws.Range("M4").Value = 0.06942796355718395
$ws.Range("N4").Value = 0.2917680520961056
